$d = $word.ActiveDocument

# Replace the placeholder text {primary-text} with {caption-text}
$d.Content.Find.Execute("{primary-text}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{caption-text}", 2)

# Remove the leftover _GoBack bookmark
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
